$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 64, pushing the existing rows 64..181 down to
# 65..182 (this also grows the sheet dimension to A1:R182 automatically).
$ws.Rows.Item(64).EntireRow.Insert()

# Populate the newly inserted row 64 with the new weekly record.
$ws.Range("A64").Value = 3
$ws.Range("B64").Value = "Femacal de La Calera"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = 44469
$ws.Range("E64").Value = 5
$ws.Range("F64").Value = 100112039
$ws.Range("G64").Value = "Ciboulette"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 160
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 1500
$ws.Range("M64").Value = 1500
$ws.Range("N64").Value = "$/docena de atados"
$ws.Range("O64").Value = "Provincia de Quillota"
$ws.Range("P64").Value = 500
$ws.Range("Q64").Value = 3
$ws.Range("R64").Value = "Hortaliza"
